$wb = $excel.ActiveWorkbook

# Sheet "All"
$wsAll = $wb.Worksheets.Item("All")
$wsAll.Range("A2").Value = "2022 September 28"
$wsAll.Range("D8").Value = "JS220_fp_usb2_ba_pcb_revB"

# Sheet "Top"
$wsTop = $wb.Worksheets.Item("Top")
$wsTop.Range("A2").Value = "2022 September 28"
$wsTop.Range("D8").Value = "JS220_fp_usb2_ba_pcb_revB"

# Sheet "Bottom"
$wsBottom = $wb.Worksheets.Item("Bottom")
$wsBottom.Range("A2").Value = "2022 September 28"
